# Weekly update: insert two new "Alcachofa" price report rows for
# Femacal de La Calera (Coquimbo) ahead of the existing history, shifting
# all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 535-536, pushing old rows 535.. down to 537..
$ws.Range("A535:A536").EntireRow.Insert()

# --- New row 535 ---
$ws.Range("A535").Value = 3
$ws.Range("B535").Value = 'Femacal de La Calera'
$ws.Range("C535").Value = 'Coquimbo'
$ws.Range("D535").Value = 45180
$ws.Range("E535").Value = 5
$ws.Range("F535").Value = 100112013
$ws.Range("G535").Value = 'Alcachofa'
$ws.Range("H535").Value = 'Española'
$ws.Range("I535").Value = 'Primera'
$ws.Range("J535").Value = 90
$ws.Range("K535").Value = 9000
$ws.Range("L535").Value = 9000
$ws.Range("M535").Value = 9000
$ws.Range("N535").Value = '$/caja 30 unidades'
$ws.Range("O535").Value = 'Provincia de Limarí'
$ws.Range("P535").Value = 300
$ws.Range("Q535").Value = 30
$ws.Range("R535").Value = 'Hortaliza'

# --- New row 536 ---
$ws.Range("A536").Value = 3
$ws.Range("B536").Value = 'Femacal de La Calera'
$ws.Range("C536").Value = 'Coquimbo'
$ws.Range("D536").Value = 45180
$ws.Range("E536").Value = 5
$ws.Range("F536").Value = 100112013
$ws.Range("G536").Value = 'Alcachofa'
$ws.Range("H536").Value = 'Española'
$ws.Range("I536").Value = 'Segunda'
$ws.Range("J536").Value = 60
$ws.Range("K536").Value = 7000
$ws.Range("L536").Value = 7000
$ws.Range("M536").Value = 7000
$ws.Range("N536").Value = '$/caja 40 unidades'
$ws.Range("O536").Value = 'Provincia de Limarí'
$ws.Range("P536").Value = 175
$ws.Range("Q536").Value = 40
$ws.Range("R536").Value = 'Hortaliza'
